$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: merge the lone-space run with the following
# "Adult females have a survival probability " run into a single run.
# ---------------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute(
    "Adult females have a survival probability", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Adult females have a survival probability", 2)
Write-Output "Hunk1 Found: $ok1"

# ---------------------------------------------------------------------------
# Hunk 2: "The two homologous locus in females act additively, and are
# averaged to obtain the phenotypic values" -> split into three runs, with
# "locus" corrected to "loci".
# ---------------------------------------------------------------------------
$t = $d.Content.Text
$idxLocus = $t.IndexOf("locus in females act additively")

# Replace locus -> loci, tightly scoped to the word itself.
$rngWord = $d.Range($idxLocus, $idxLocus + 5)
$ok2 = $rngWord.Find.Execute("locus", $true, $false, $false, $false, $false,
                              $true, 1, $false, "loci", 2)
Write-Output "Hunk2 Found: $ok2"

# Recompute offsets against the post-replace text (it is one character
# shorter than before).
$t2 = $d.Content.Text
$idxTwo = $t2.IndexOf("The two homologous loci")
$idxLociStart = $idxTwo + ("The two homologous ").Length
$idxLociEnd = $idxLociStart + ("loci").Length
$idxValuesEnd = $t2.IndexOf(" << Again confirm", $idxTwo)

# Re-establish the run boundary between "Genes are expressed..." and
# "The two homologous ..." (a plain formatting no-op toggle splits runs
# without re-merging across the paragraph).
$rngA = $d.Range($idxTwo, $idxValuesEnd)
$rngA.Font.Bold = 1
$rngA.Font.Bold = 0

# Split "The two homologous " from "loci".
$rngB = $d.Range($idxLociStart, $idxValuesEnd)
$rngB.Font.Bold = 1
$rngB.Font.Bold = 0

# Split "loci" from " in females act additively, and are averaged to obtain
# the phenotypic values".
$rngC = $d.Range($idxLociEnd, $idxValuesEnd)
$rngC.Font.Bold = 1
$rngC.Font.Bold = 0

Write-Output "Done"
